# Applies the tracked changes described by the commit's unified diff to
# IPC/questões das boas.docx. Uses Range.InsertXML for the edits that split
# an existing run into several runs (so the resulting OOXML run layout
# matches the diff exactly), and plain Range/Find based text edits for the
# simple in-place word changes.

$d = $word.ActiveDocument

function Get-DocText {
    return $d.Content.Text
}

function Replace-ByOffset($needle, $xmlBody) {
    # Locates the first occurrence of $needle in the document and replaces
    # that exact character range with the supplied OOXML body (a sequence
    # of <w:r>/<w:bookmarkStart>/<w:bookmarkEnd> etc. elements that belongs
    # inside a <w:p>).
    # NOTE: positional parameters only -- this interpreter does not bind
    # PowerShell-style "-name value" named arguments correctly.
    $full = Get-DocText
    $idx = $full.IndexOf($needle)
    if ($idx -lt 0) {
        throw "Could not locate text: $needle"
    }
    $target = $d.Range($idx, $idx + $needle.Length)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $xmlBody + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the title paragraph down to the end
#    of the "utilizadores do sistema ... Beja." paragraph, and change /
#    split that paragraph's trailing sentence.
# ---------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rPrCambria = '<w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/></w:rPr>'

Replace-ByOffset " constituídos por alunos e docentes do Instituto Politécnico de Beja." (
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve"> constit</w:t></w:r>' +
    '<w:r>' + $rPrCambria + '<w:t>uídos por alunos e docentes da Escola Superior de Tecnologia e Gestão do Instituto Politécnico de Beja.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
)

# ---------------------------------------------------------------------
# 2) "Consultar a sua assiduidade ..." -> "Consultar a assiduidade ..."
#    (simple in-run word removal, no run split in the diff)
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    "Consultar a sua assiduidade numa determinada unidade curricular.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Consultar a assiduidade numa determinada unidade curricular.", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) "Criar uma aula para a uma determinada unidade curricular;" is split
#    into three runs and "a " is dropped ("para a uma" -> "para uma").
# ---------------------------------------------------------------------

Replace-ByOffset "Criar uma aula para a uma determinada unidade curricular;" (
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve">Criar </w:t></w:r>' +
    '<w:r>' + $rPrCambria + '<w:t>uma aula para</w:t></w:r>' +
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve"> uma determinada unidade curricular;</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 4) "...se familiarizarem com o novo sistema, de forma absolutamente
#    intuitiva." -> split into three runs, drop the comma and
#    "absolutamente ".
# ---------------------------------------------------------------------

Replace-ByOffset ") no processo de desenho da interface, que irão permitir aos utilizadores realizarem uma reciclagem de conhecimentos, bem como se familiarizarem com o novo sistema, de forma absolutamente intuitiva." (
    '<w:r>' + $rPrCambria + '<w:t>) no processo de desenho da interface, que irão permitir aos utilizadores realizarem uma reciclagem de conhecimentos, bem como se fam</w:t></w:r>' +
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve">iliarizarem com o novo sistema </w:t></w:r>' +
    '<w:r>' + $rPrCambria + '<w:t>de forma intuitiva.</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 5) "...de modo a poder consultar..." -> "...de modo a ser possível
#    consultar..." split into three runs.
# ---------------------------------------------------------------------

Replace-ByOffset ". Todos eles necessitam de uma conexão válida à Internet e da aplicação instada no seu dispositivo, de modo a poder consultar a informação disponibilizada no sistema." (
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve">. Todos eles necessitam de uma conexão válida à Internet e da aplicação instada no seu dispositivo, de modo a </w:t></w:r>' +
    '<w:r>' + $rPrCambria + '<w:t>ser possível</w:t></w:r>' +
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve"> consultar a informação disponibilizada no sistema.</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 6) "O sistema não deverá possuir..." -> "A aplicação não deverá
#    possuir..." split into two runs.
# ---------------------------------------------------------------------

Replace-ByOffset "O sistema não deverá possuir um sistema de comunicação, visto que detém um carácter informativo e de inserção de dados, dispensando a necessidade de implementação de uma funcionalidade de comunicação entre utilizadores." (
    '<w:r>' + $rPrCambria + '<w:t xml:space="preserve">A aplicação </w:t></w:r>' +
    '<w:r>' + $rPrCambria + '<w:t>não deverá possuir um sistema de comunicação, visto que detém um carácter informativo e de inserção de dados, dispensando a necessidade de implementação de uma funcionalidade de comunicação entre utilizadores.</w:t></w:r>'
)

Write-Host "All edits applied."
